$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("10" - 10kΩ resistor group): R37 removed from the group (moved to new 0Ω group)
$ws.Range("B11").Value = 17
$ws.Range("D11").Value = "R1,R3,R4,R5,R6,R18,R19,R20,R21,R22,R23,R24,R25,R33,R34,R35,R36"

# Row 15 ("14"): SW1 switch entry replaced with a new 0Ω resistor entry (R37,R38,R39)
$ws.Range("B15").Value = 3
$ws.Range("C15").Value = "0Ω"
$ws.Range("D15").Value = "R37,R38,R39"
$ws.Range("E15").Value = "R0805"
$ws.Range("F15").Value = "0Ω"
$ws.Range("G15").Value = "0805W8F0000T5E"
$ws.Range("H15").Value = "UNI-ROYAL(厚声)"
$ws.Range("I15").Value = "C17477"
$ws.Range("J15").Value = "LCSC"
